# Daily attendance processing - 2026-02-01 13:57:10
# For every row in the "Recorded By" column (G), swap the order of the
# first two comma-separated names/emails (any trailing extra entries,
# e.g. a third "system" marker, are left in place at the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $text = [string]$cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ', '
    if ($parts.Count -ge 2) {
        $swapped = @($parts[1], $parts[0])
        if ($parts.Count -gt 2) {
            $swapped += $parts[2..($parts.Count - 1)]
        }
        $newText = [string]::Join(', ', $swapped)
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
